$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.677.07"
$ws.Range("E2").Value = "  +2.16%  "

$ws.Range("D3").Value = "2.157.87"
$ws.Range("E3").Value = "  +2.38%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.63"
$ws.Range("E5").Value = "  -0.42%  "

$ws.Range("E6").Value = "  +2.29%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.50"
$ws.Range("E7").Value = "  +1.74%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.392"
$ws.Range("E9").Value = "  +0.83%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0845"
$ws.Range("E10").Value = "  +0.21%  "

$ws.Range("E11").Value = "  +0.08%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.90"
$ws.Range("E12").Value = "  +0.37%  "

$ws.Range("D13").Value = "2.481.23"
$ws.Range("E13").Value = "  +2.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.89"
$ws.Range("E14").Value = "  -0.92%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.806"
$ws.Range("E15").Value = "  -0.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.48"
$ws.Range("E16").Value = "  -0.96%  "

$ws.Range("D17").Value = "2.162.75"
$ws.Range("E17").Value = "  +1.74%  "

$ws.Range("D18").Value = "39.602.77"
$ws.Range("E18").Value = "  +1.93%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.80"
$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.08"
$ws.Range("E20").Value = "  -0.47%  "

$ws.Range("D21").Value = "0.0₃0843"
$ws.Range("E21").Value = "  -0.51%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "231.07"
$ws.Range("E22").Value = "  +1.29%  "

$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.39"
$ws.Range("E24").Value = "  +3.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").Value = "  -2.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "172.83"
$ws.Range("E26").Value = "  +0.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.60"
$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("E28").Value = "  +0.86%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.98"
$ws.Range("E29").Value = "  +3.16%  "

$ws.Range("E30").Value = "  +0.28%  "

$ws.Range("E31").Value = "  +4.72%  "

$ws.Range("E32").Value = "  +1.22%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.57"
$ws.Range("E33").Value = "  -0.16%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.66"
$ws.Range("E34").Value = "  -1.91%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.94"
$ws.Range("E35").Value = "  -3.65%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0617"
$ws.Range("E36").Value = "  -0.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.39"
$ws.Range("E37").Value = "  -0.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.59"
$ws.Range("E38").Value = "  +1.87%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.14"
$ws.Range("E39").Value = "  +23.59%  "

$ws.Range("E40").Value = "  -0.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.18"
$ws.Range("E41").Value = "  +0.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0226"
$ws.Range("E42").Value = "  -0.37%  "

$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.40"
$ws.Range("E43").Value = "  -3.64%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.511.21"
$ws.Range("E44").Value = "  -0.92%  "

$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.22"
$ws.Range("E45").Value = "  +1.22%  "

$ws.Range("B46").Value = "HuobiToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.81"
$ws.Range("E46").Value = "  +0.26%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0919"
$ws.Range("E47").Value = "  +0.62%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.09"
$ws.Range("E48").Value = "  +0.29%  "

$ws.Range("E49").Value = "  -0.67%  "

$ws.Range("D51").Value = "2.364.36"
$ws.Range("E51").Value = "  +2.61%  "
